# Runmode mode functionality added

$wb = $excel.ActiveWorkbook

# Remove the LoginTest and CreateAccountTest sheets, keep only the first sheet.
$wb.Worksheets.Item("LoginTest").Delete() | Out-Null
$wb.Worksheets.Item("CreateAccountTest").Delete() | Out-Null

$ws = $wb.Worksheets.Item("test_suite")
$ws.Name = "findNewCar"

# Rewrite the data with the new browser/runmode layout.
$ws.Range("A1").Value = "browser"
$ws.Range("B1").Value = "runmode"
$ws.Range("A2").Value = "chrome"
$ws.Range("B2").Value = "Y"
$ws.Range("A3").Value = "chrome"
$ws.Range("B3").Value = "N"
$ws.Range("A4").Value = "chrome"
$ws.Range("B4").Value = "Y"

$ws.Columns.Item(2).ColumnWidth = 20.1

$ws.Range("B4").Select() | Out-Null
